# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# values on Sheet1, matching the automated "Updated symbol list" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    # A leading apostrophe forces Excel to store the value as literal text
    # (preserving things like "273.32" and "-1.94%" exactly, without being
    # re-interpreted as a number/percentage).
    $c.Value = "'" + $text
    # Drop the quote-prefix style Excel applies for text-forced entries so
    # the cell keeps its original (unstyled) appearance.
    $c.Style = "Normal"
}

# Row -> (Price, Volume%) new values, taken from the diff.
$updates = @(
    @{ Row = 2;  D = "273.32";     E = "-1.94%" }
    @{ Row = 3;  D = "26.58";      E = "-2.55%" }
    @{ Row = 4;  D = "4.883";      E = "1.85%" }
    @{ Row = 5;  D = $null;        E = "1.32%" }
    @{ Row = 6;  D = "6.906";      E = "0.81%" }
    @{ Row = 7;  D = "3.358";      E = "2.74%" }
    @{ Row = 8;  D = "1.257";      E = "33.90%" }
    @{ Row = 9;  D = "0.8774";     E = "0.16%" }
    @{ Row = 10; D = "0.1453";     E = "0.19%" }
    @{ Row = 11; D = "0.05138";    E = "1.69%" }
    @{ Row = 12; D = "0.07325";    E = "0.66%" }
    @{ Row = 13; D = "0.03125";    E = "-0.71%" }
    @{ Row = 14; D = "0.09045";    E = "0.11%" }
    @{ Row = 15; D = "0.001577";   E = "2.23%" }
    @{ Row = 16; D = "0.0006318";  E = "0.93%" }
    @{ Row = 17; D = "0.006044";   E = "-0.15%" }
    @{ Row = 18; D = "3.455";      E = "-0.37%" }
    @{ Row = 19; D = $null;        E = "-0.11%" }
    @{ Row = 21; D = "0.1327";     E = "1.35%" }
    @{ Row = 22; D = "3.906";      E = "1.28%" }
    @{ Row = 23; D = "0.04428";    E = "2.51%" }
    @{ Row = 24; D = "0.001178";   E = "0.32%" }
    @{ Row = 25; D = $null;        E = "3.16%" }
    @{ Row = 27; D = $null;        E = "5.61%" }
    @{ Row = 40; D = "0.04030";    E = "0.05%" }
    @{ Row = 41; D = "0.006660";   E = "-0.62%" }
    @{ Row = 42; D = $null;        E = "1.35%" }
    @{ Row = 43; D = "0.002103";   E = "-1.13%" }
    @{ Row = 44; D = "0.01252";    E = "-6.46%" }
    @{ Row = 45; D = "0.00005319"; E = "3.83%" }
    @{ Row = 46; D = $null;        E = "8.16%" }
    @{ Row = 47; D = "0.02003";    E = "-32.95%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextCell $u.Row 4 $u.D
    }
    if ($null -ne $u.E) {
        Set-TextCell $u.Row 5 $u.E
    }
}
